$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.834.55"
$ws.Range("E2").Value = "  +0.73%  "

$ws.Range("D3").Value = "2.639.70"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.60"
$ws.Range("E5").Value = "  -1.31%  "

$ws.Range("D6").Value = "146.99"
$ws.Range("E6").Value = "  +2.12%  "

$ws.Range("E8").Value = "  +0.38%  "

$ws.Range("E9").Value = "  +1.84%  "

$ws.Range("D10").Value = "0.384"
$ws.Range("E10").Value = "  +6.21%  "

$ws.Range("D11").Value = "5.59"
$ws.Range("E11").Value = "  -0.48%  "

$ws.Range("E12").Value = "  -1.05%  "

$ws.Range("D13").Value = "27.46"
$ws.Range("E13").Value = "  +0.15%  "

$ws.Range("D14").Value = "3.112.69"
$ws.Range("E14").Value = "  -1.43%  "

$ws.Range("D15").Value = "63.667.70"
$ws.Range("E15").Value = "  +0.67%  "

$ws.Range("E16").Value = "  +1.39%  "

$ws.Range("D17").Value = "2.634.34"
$ws.Range("E17").Value = "  -1.64%  "

$ws.Range("D18").Value = "11.75"
$ws.Range("E18").Value = "  +2.60%  "

$ws.Range("D19").Value = "4.56"
$ws.Range("E19").Value = "  +3.47%  "

$ws.Range("D20").Value = "346.94"
$ws.Range("E20").Value = "  +1.22%  "

$ws.Range("D21").Value = "6.92"
$ws.Range("E21").Value = "  +0.76%  "

$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("E23").Value = "  -1.23%  "

$ws.Range("D24").Value = "66.31"
$ws.Range("E24").Value = "  -1.46%  "

$ws.Range("D25").Value = "1.67"
$ws.Range("E25").Value = "  +8.04%  "

$ws.Range("D26").Value = "1.69"
$ws.Range("E26").Value = "  +2.30%  "

$ws.Range("D27").Value = "9.24"
$ws.Range("E27").Value = "  +6.43%  "

$ws.Range("D28").Value = "564.74"
$ws.Range("E28").Value = "  +5.01%  "

$ws.Range("D29").Value = "8.11"
$ws.Range("E29").Value = "  +2.40%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.03%  "

$ws.Range("E31").Value = "  -2.96%  "

$ws.Range("E32").Value = "  -1.11%  "

$ws.Range("D33").Value = "0.0₃0851"
$ws.Range("E33").Value = "  +5.20%  "

$ws.Range("E34").Value = "  -1.78%  "

$ws.Range("E35").Value = "  +1.98%  "

$ws.Range("D36").Value = "169.05"
$ws.Range("E36").Value = "  -1.88%  "

$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.08%  "

$ws.Range("D38").Value = "0.405"
$ws.Range("E38").Value = "  -0.27%  "

$ws.Range("D39").Value = "1.95"
$ws.Range("E39").Value = "  +4.43%  "

$ws.Range("D40").Value = "19.18"
$ws.Range("E40").Value = "  -0.37%  "

$ws.Range("E41").Value = "  +0.00%  "

$ws.Range("D42").Value = "165.16"
$ws.Range("E42").Value = "  -6.69%  "

$ws.Range("D43").Value = "40.07"
$ws.Range("E43").Value = "  -0.13%  "

$ws.Range("D44").Value = "3.79"
$ws.Range("E44").Value = "  +0.98%  "

$ws.Range("D45").Value = "21.95"
$ws.Range("E45").Value = "  -1.56%  "

$ws.Range("D46").Value = "0.0567"
$ws.Range("E46").Value = "  -0.78%  "

$ws.Range("D47").Value = "0.626"
$ws.Range("E47").Value = "  -1.63%  "

$ws.Range("E48").Value = "  +2.30%  "

$ws.Range("E49").Value = "  +13.72%  "

$ws.Range("D50").Value = "0.0957"
$ws.Range("E50").Value = "  -0.83%  "

$ws.Range("E51").Value = "  -0.37%  "
